$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "44.093.02"
$ws.Range("E2").Value = "  +1.63%  "

# Row 3
$ws.Range("D3").Value = "2.273.18"
$ws.Range("E3").Value = "  +0.47%  "

# Row 4
$ws.Range("E4").Value = "  -0.48%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.52"
$ws.Range("E5").Value = "  +0.16%  "

# Row 6
$ws.Range("E6").Value = "  +1.57%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.11"
$ws.Range("E7").Value = "  +5.05%  "

# Row 8
$ws.Range("E8").Value = "  -0.10%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.450"
$ws.Range("E9").Value = "  +9.61%  "

# Row 10
$ws.Range("E10").Value = "  +10.25%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.31"
$ws.Range("E11").Value = "  -0.16%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.92"
$ws.Range("E12").Value = "  +19.70%  "

# Row 13
$ws.Range("E13").Value = "  +1.76%  "

# Row 14
$ws.Range("D14").Value = "2.613.34"
$ws.Range("E14").Value = "  +0.32%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.72"
$ws.Range("E15").Value = "  +0.46%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.14"
$ws.Range("E16").Value = "  +8.18%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.842"
$ws.Range("E17").Value = "  +4.40%  "

# Row 18
$ws.Range("D18").Value = "2.274.42"
$ws.Range("E18").Value = "  +0.22%  "

# Row 19
$ws.Range("D19").Value = "44.017.05"
$ws.Range("E19").Value = "  +1.62%  "

# Row 20
$ws.Range("E20").Value = "  +8.74%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.91"
$ws.Range("E21").Value = "  +1.35%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.11"
$ws.Range("E22").Value = "  -1.25%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.01"
$ws.Range("E23").Value = "  +2.53%  "

# Row 24
$ws.Range("E24").Value = "  +0.00%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("E25").Value = "  -4.40%  "

# Row 26
$ws.Range("E26").Value = "  -3.07%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.37"
$ws.Range("E27").Value = "  +25.61%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  +3.15%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.64"
$ws.Range("E29").Value = "  +1.25%  "

# Row 30
$ws.Range("E30").Value = "  -1.82%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.86"
$ws.Range("E31").Value = "  +1.62%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.40"
$ws.Range("E32").Value = "  -5.22%  "

# Row 33
$ws.Range("E33").Value = "  +3.37%  "

# Row 34
$ws.Range("E34").Value = "  +7.34%  "

# Row 35
$ws.Range("E35").Value = "  +2.42%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.92"
$ws.Range("E36").Value = "  -1.75%  "

# Row 37
$ws.Range("E37").Value = "  +5.89%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.55"
$ws.Range("E38").Value = "  +2.04%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.33"
$ws.Range("E39").Value = "  -2.42%  "

# Row 40
$ws.Range("E40").Value = "  +4.02%  "

# Row 41
$ws.Range("E41").Value = "  +10.64%  "

# Row 42
$ws.Range("E42").Value = "  -0.22%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.61"
$ws.Range("E43").Value = "  +5.51%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0979"
$ws.Range("E44").Value = "  +1.03%  "

# Row 45
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.54"
$ws.Range("E45").Value = "  +20.07%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.26"
$ws.Range("E46").Value = "  -4.50%  "

# Row 47
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.47"
$ws.Range("E47").Value = "  +2.01%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.34"
$ws.Range("E48").Value = "  +1.70%  "

# Row 49
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.21"
$ws.Range("E49").Value = "  +0.19%  "

# Row 50
$ws.Range("D50").Value = "1.449.71"
$ws.Range("E50").Value = "  -0.79%  "

# Row 51
$ws.Range("E51").Value = "  +3.83%  "
